$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 0.005723000769734084
$ws.Range("J2").Value = 0.005723000769734084
$ws.Range("O2").Value = 0.0002880357555630755
$ws.Range("P2").Value = 0.0002880357555630755
$ws.Range("S2").Value = 0.00000164842885079842
$ws.Range("T2").Value = 0.00000164842885079842
$ws.Range("I3").Value = 0.005723000769734084
$ws.Range("J3").Value = 0.005723000769734084
$ws.Range("M3").Value = 32.87103466666667
$ws.Range("N3").Value = 98.61310400000001
$ws.Range("O3").Value = 0.7767049471988007
$ws.Range("P3").Value = 0.7767049471988008
$ws.Range("Q3").Value = 9.274715829361778
$ws.Range("R3").Value = 83.47244246425601
$ws.Range("S3").Value = 0.004445083010675008
$ws.Range("T3").Value = 0.004445083010675008
$ws.Range("I4").Value = 0.005723000769734084
$ws.Range("J4").Value = 0.005723000769734084
$ws.Range("M4").Value = 9.43791
$ws.Range("N4").Value = 28.31373
$ws.Range("O4").Value = 0.2230070170456362
$ws.Range("P4").Value = 0.2230070170456362
$ws.Range("Q4").Value = 2.66295035008
$ws.Range("R4").Value = 23.96655315072
$ws.Range("S4").Value = 0.001276269330208278
$ws.Range("T4").Value = 0.001276269330208278
$ws.Range("G5").Value = 30.199365
$ws.Range("H5").Value = 90.598095
$ws.Range("I5").Value = 0.6125398923302606
$ws.Range("J5").Value = 0.6125398923302606
$ws.Range("O5").Value = 0.0002880357555630755
$ws.Range("P5").Value = 0.0002880357555630755
$ws.Range("Q5").Value = 0.36813025935
$ws.Range("R5").Value = 3.31317233415
$ws.Range("S5").Value = 0.0001764333906998715
$ws.Range("T5").Value = 0.0001764333906998715
$ws.Range("G6").Value = 30.199365
$ws.Range("H6").Value = 90.598095
$ws.Range("I6").Value = 0.6125398923302606
$ws.Range("J6").Value = 0.6125398923302606
$ws.Range("M6").Value = 32.87103466666667
$ws.Range("N6").Value = 98.61310400000001
$ws.Range("O6").Value = 0.7767049471988007
$ws.Range("P6").Value = 0.7767049471988008
$ws.Range("Q6").Value = 992.68437382632
$ws.Range("R6").Value = 8934.15936443688
$ws.Range("S6").Value = 0.4757627647295341
$ws.Range("T6").Value = 0.4757627647295342
$ws.Range("G7").Value = 30.199365
$ws.Range("H7").Value = 90.598095
$ws.Range("I7").Value = 0.6125398923302606
$ws.Range("J7").Value = 0.6125398923302606
$ws.Range("M7").Value = 9.43791
$ws.Range("N7").Value = 28.31373
$ws.Range("O7").Value = 0.2230070170456362
$ws.Range("P7").Value = 0.2230070170456362
$ws.Range("Q7").Value = 285.01888892715
$ws.Range("R7").Value = 2565.17000034435
$ws.Range("S7").Value = 0.1366006942100266
$ws.Range("T7").Value = 0.1366006942100266
$ws.Range("G8").Value = 18.820355
$ws.Range("H8").Value = 56.461065
$ws.Range("I8").Value = 0.3817371069000054
$ws.Range("J8").Value = 0.3817371069000054
$ws.Range("O8").Value = 0.0002880357555630755
$ws.Range("P8").Value = 0.0002880357555630755
$ws.Range("Q8").Value = 0.22942012745
$ws.Range("R8").Value = 2.06478114705
$ws.Range("S8").Value = 0.0001099539360124056
$ws.Range("T8").Value = 0.0001099539360124056
$ws.Range("G9").Value = 18.820355
$ws.Range("H9").Value = 56.461065
$ws.Range("I9").Value = 0.3817371069000054
$ws.Range("J9").Value = 0.3817371069000054
$ws.Range("M9").Value = 32.87103466666667
$ws.Range("N9").Value = 98.61310400000001
$ws.Range("O9").Value = 0.7767049471988007
$ws.Range("P9").Value = 0.7767049471988008
$ws.Range("Q9").Value = 618.6445416439733
$ws.Range("R9").Value = 5567.80087479576
$ws.Range("S9").Value = 0.2964970994585916
$ws.Range("T9").Value = 0.2964970994585917
$ws.Range("G10").Value = 18.820355
$ws.Range("H10").Value = 56.461065
$ws.Range("I10").Value = 0.3817371069000054
$ws.Range("J10").Value = 0.3817371069000054
$ws.Range("M10").Value = 9.43791
$ws.Range("N10").Value = 28.31373
$ws.Range("O10").Value = 0.2230070170456362
$ws.Range("P10").Value = 0.2230070170456362
$ws.Range("Q10").Value = 177.62481665805
$ws.Range("R10").Value = 1598.62334992245
$ws.Range("S10").Value = 0.08513005350540137
$ws.Range("T10").Value = 0.08513005350540137
